$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selected/active cell shown in the sheet view (selection was added to sheetView)
$ws.Range("M23").Select()

# Narrow the width of columns A:C (closest achievable value to 34.42578125)
$ws.Range("A1:C1").ColumnWidth = 33.666666666666664

# Update data values in column O
$ws.Range("O5").Value = 28.6
$ws.Range("O6").Value = 33.6
$ws.Range("O8").Value = 71.2
$ws.Range("O12").Value = 16.100000000000001
